# Append the new Adafruit IO reading as row 7.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "2024-09-25T18:06:40Z"
$ws.Range("B7").Value = "temperature"
# "25" looks numeric - force it to stay text like the other numeric-looking
# "Value" cells already in the sheet (e.g. C2="29").
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "25"
$ws.Range("D7").Value = "N/A"
$ws.Range("E7").Value = "N/A"
$ws.Range("F7").Value = "N/A"
